$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 'isophonics_41'
$ws.Range("B2").Value = 'isophonics_288'
$ws.Range("C2").Value = 0.06944444444444445
$ws.Range("D2").Value = '[[''G'', ''C'', ''F'']]'
$ws.Range("E2").Value = '[[''F#'', ''B'', ''E'']]'
$ws.Range("F2").Value = '[(56.279931, 60.633673)]'
$ws.Range("G2").Value = '[(29.008095, 34.615714)]'
$ws.Range("H2").Value = ""
$ws.Range("I2").Value = ""

# Row 3
$ws.Range("A3").Value = 'schubert-winterreise_163'
$ws.Range("B3").Value = 'schubert-winterreise_145'
$ws.Range("C3").Value = 0.5333333333333333
$ws.Range("D3").Value = '[[''F:maj'', ''A#:maj'', ''F:maj'', ''A#:maj'']]'
$ws.Range("E3").Value = '[[''D:maj/A'', ''G:maj'', ''D:maj/A'', ''G:maj/B'']]'
$ws.Range("F3").Value = '[(111.92, 115.66)]'
$ws.Range("G3").Value = '[(143.58, 149.12)]'
$ws.Range("H3").Value = 'spotify:track:1nvxQGWCnikMK7a4HYQvSx'
$ws.Range("I3").Value = ""

# Row 4
$ws.Range("A4").Value = 'jaah_0'
$ws.Range("B4").Value = 'schubert-winterreise_200'
$ws.Range("C4").Value = 0.3189655172413793
$ws.Range("D4").Value = '[[''Ab'', ''Eb:7'', ''Ab'', ''Ab'']]'
$ws.Range("E4").Value = '[[''E:maj'', ''B:7/A'', ''E:maj/G#'', ''E:maj/B'']]'
$ws.Range("F4").Value = '[(63.01, 68.27)]'
$ws.Range("G4").Value = '[(207.92, 219.14)]'
$ws.Range("H4").Value = ""
$ws.Range("I4").Value = ""

# Row 5
$ws.Range("A5").Value = 'schubert-winterreise_163'
$ws.Range("B5").Value = 'schubert-winterreise_2'
$ws.Range("C5").Value = 0.3939393939393939
$ws.Range("D5").Value = '[[''F:maj'', ''C:7'', ''F:maj'', ''C:7'', ''F:maj'']]'
$ws.Range("E5").Value = '[[''A:maj/E'', ''E:7'', ''A:maj'', ''E:7'', ''A:maj'']]'
$ws.Range("F5").Value = '[(55.4, 65.22)]'
$ws.Range("G5").Value = '[(20.56, 26.4)]'
$ws.Range("H5").Value = 'spotify:track:1nvxQGWCnikMK7a4HYQvSx'
$ws.Range("I5").Value = 'spotify:track:0XfunCHFEeQnzm4NaY8rJr'

# Row 6
$ws.Range("A6").Value = 'schubert-winterreise_44'
$ws.Range("B6").Value = 'schubert-winterreise_214'
$ws.Range("C6").Value = 0.1916666666666667
$ws.Range("D6").Value = '[[''A#/F'', ''F:7'', ''A#'']]'
$ws.Range("E6").Value = '[[''G:maj'', ''D:7'', ''G:maj'']]'
$ws.Range("F6").Value = '[(271.22, 275.78)]'
$ws.Range("G6").Value = '[(9.2, 20.1)]'
$ws.Range("H6").Value = ""
$ws.Range("I6").Value = ""

# Row 7
$ws.Range("A7").Value = 'schubert-winterreise_149'
$ws.Range("B7").Value = 'schubert-winterreise_21'
$ws.Range("C7").Value = 0.150268336314848
$ws.Range("D7").Value = '[[''C:maj/G'', ''G:7'', ''C:maj'']]'
$ws.Range("E7").Value = '[[''F#:maj'', ''C#:7/F'', ''F#:maj'']]'
$ws.Range("F7").Value = '[(32.44, 34.8)]'
$ws.Range("G7").Value = '[(38.58, 43.3)]'
$ws.Range("H7").Value = 'spotify:track:2qCvEz2hEb92VFATqVvrht'
$ws.Range("I7").Value = ""

# Row 8
$ws.Range("A8").Value = 'isophonics_212'
$ws.Range("B8").Value = 'jaah_55'
$ws.Range("C8").Value = 0.09545454545454546
$ws.Range("D8").Value = '[[''G'', ''D'', ''D'']]'
$ws.Range("E8").Value = '[[''F'', ''C'', ''C'']]'
$ws.Range("F8").Value = '[(57.636679, 62.11813)]'
$ws.Range("G8").Value = '[(37.61, 38.9)]'
$ws.Range("H8").Value = ""
$ws.Range("I8").Value = ""

# Row 9
$ws.Range("A9").Value = 'isophonics_124'
$ws.Range("B9").Value = 'isophonics_167'
$ws.Range("C9").Value = 0.162280701754386
$ws.Range("D9").Value = '[[''F'', ''Bb'', ''F'']]'
$ws.Range("E9").Value = '[[''G/5'', ''C/9'', ''G'']]'
$ws.Range("F9").Value = '[(15.124263, 20.058503)]'
$ws.Range("G9").Value = '[(0.727609, 8.204625)]'
$ws.Range("H9").Value = ""
$ws.Range("I9").Value = ""

# Row 10
$ws.Range("A10").Value = 'isophonics_1'
$ws.Range("B10").Value = 'isophonics_140'
$ws.Range("C10").Value = 0.202991452991453
$ws.Range("D10").Value = '[[''Eb'', ''Ab/5'', ''Eb'', ''Ab/5'', ''Eb'']]'
$ws.Range("E10").Value = '[[''G'', ''C'', ''G'', ''C'', ''G'']]'
$ws.Range("F10").Value = '[(17.016, 25.693)]'
$ws.Range("G10").Value = '[(0.465952, 8.50033)]'
$ws.Range("H10").Value = ""
$ws.Range("I10").Value = 'spotify:track:3VbGCXWRiouAq8VyMYN2MI'

# Row 11
$ws.Range("A11").Value = 'schubert-winterreise_34'
$ws.Range("B11").Value = 'schubert-winterreise_154'
$ws.Range("C11").Value = 0.1613636363636363
$ws.Range("D11").Value = '[[''C/G'', ''G:7'', ''C''], [''C'', ''C/G'', ''G:7'']]'
$ws.Range("E11").Value = '[[''A:maj/E'', ''E:7'', ''A:maj''], [''A:maj'', ''A:maj'', ''E:7/G#'']]'
$ws.Range("F11").Value = '[(243.44, 246.98), (242.22, 244.56)]'
$ws.Range("G11").Value = '[(21.78, 25.3), (10.3, 17.72)]'
$ws.Range("H11").Value = ""
$ws.Range("I11").Value = 'spotify:track:0XfunCHFEeQnzm4NaY8rJr'

# Row 12
$ws.Range("A12").Value = 'schubert-winterreise_78'
$ws.Range("B12").Value = 'schubert-winterreise_21'
$ws.Range("C12").Value = 0.3311965811965812
$ws.Range("D12").Value = '[[''A:maj'', ''D:min'', ''A:maj'', ''D:min'', ''A:maj''], [''D:min'', ''A:maj'', ''D:min'', ''A:maj'', ''D:min'']]'
$ws.Range("E12").Value = '[[''A#:maj'', ''D#:min'', ''A#:maj'', ''D#:min'', ''A#:maj''], [''D#:min'', ''A#:maj'', ''D#:min'', ''A#:maj'', ''D#:min'']]'
$ws.Range("F12").Value = '[(6.22, 14.32), (7.4, 16.48)]'
$ws.Range("G12").Value = '[(2.02, 8.58), (2.66, 10.62)]'
$ws.Range("H12").Value = ""
$ws.Range("I12").Value = ""

# Row 13
$ws.Range("A13").Value = 'schubert-winterreise_67'
$ws.Range("B13").Value = 'schubert-winterreise_26'
$ws.Range("C13").Value = 0.2363636363636364
$ws.Range("D13").Value = '[[''B:min'', ''F#:7/C#'', ''B:min/D''], [''B:min'', ''F#:maj/A#'', ''B:min'']]'
$ws.Range("E13").Value = '[[''F:min'', ''C:7'', ''F:min''], [''F:min'', ''C:maj'', ''F:min'']]'
$ws.Range("F13").Value = '[(0.3, 2.82), (16.3, 19.52)]'
$ws.Range("G13").Value = '[(9.62, 14.18), (1.88, 6.56)]'
$ws.Range("H13").Value = ""
$ws.Range("I13").Value = 'spotify:track:1nvxQGWCnikMK7a4HYQvSx'

# Row 14
$ws.Range("A14").Value = 'jaah_43'
$ws.Range("B14").Value = 'isophonics_21'
$ws.Range("C14").Value = 0.134575569358178
$ws.Range("D14").Value = '[[''Bb:7'', ''Eb'', ''Eb'', ''Ab'']]'
$ws.Range("E14").Value = '[[''G:7'', ''C'', ''C/b7'', ''F'']]'
$ws.Range("F14").Value = '[(93.64, 102.97)]'
$ws.Range("G14").Value = '[(34.041, 40.124)]'
$ws.Range("H14").Value = ""
$ws.Range("I14").Value = ""

# Row 15
$ws.Range("A15").Value = 'schubert-winterreise_27'
$ws.Range("B15").Value = 'jaah_21'
$ws.Range("C15").Value = 0.2571428571428571
$ws.Range("D15").Value = '[[''G:maj/D'', ''D:7'', ''G:maj'']]'
$ws.Range("E15").Value = '[[''Ab'', ''Eb:7'', ''Ab'']]'
$ws.Range("F15").Value = '[(64.3, 65.72)]'
$ws.Range("G15").Value = '[(46.75, 47.53)]'
$ws.Range("H15").Value = ""
$ws.Range("I15").Value = ""

# Row 16
$ws.Range("A16").Value = 'schubert-winterreise_166'
$ws.Range("B16").Value = 'schubert-winterreise_151'
$ws.Range("C16").Value = 0.3833333333333333
$ws.Range("D16").Value = '[[''D:maj'', ''G:maj'', ''D:maj/F#'']]'
$ws.Range("E16").Value = '[[''C:maj/G'', ''F:maj'', ''C:maj/G'']]'
$ws.Range("F16").Value = '[(67.14, 70.82)]'
$ws.Range("G16").Value = '[(117.54, 121.8)]'
$ws.Range("H16").Value = ""
$ws.Range("I16").Value = ""

# Row 17
$ws.Range("A17").Value = 'isophonics_227'
$ws.Range("B17").Value = 'schubert-winterreise_79'
$ws.Range("C17").Value = 0.13125
$ws.Range("D17").Value = '[[''E:7'', ''A:min'', ''A:min/b3''], [''E:7'', ''A:min'', ''D:min'']]'
$ws.Range("E17").Value = '[[''E:7'', ''A:min'', ''A:min''], [''E:7'', ''A:min'', ''D:min'']]'
$ws.Range("F17").Value = '[(5.67204, 8.841564), (20.045192, 25.141972)]'
$ws.Range("G17").Value = '[(10.24, 21.02), (53.28, 62.86)]'
$ws.Range("H17").Value = ""
$ws.Range("I17").Value = 'spotify:track:3OD2uwEUQKg0WyW9Lewata'
